# Client Setup Admin Suite - integrate client admin files
$wb = $excel.ActiveWorkbook

# --- Sheet "Test Cases" ---
$ws1 = $wb.Worksheets.Item("Test Cases")
$ws1.Range("D3").Value = "PASS"
$ws1.Range("D4").Value = "PASS"
$ws1.Range("D5").Value = "PASS"

# --- Sheet "ClientSupplierCreation" ---
$ws2 = $wb.Worksheets.Item("ClientSupplierCreation")
$ws2.Range("B3").Value = "Ajay Supplier"
$ws2.Range("D3").Value = "ajay"
$ws2.Range("K2").Value = "SKIP"
$ws2.Range("K3").Value = "PASS"

# --- Sheet "ApplicationProvisioning" ---
$ws3 = $wb.Worksheets.Item("ApplicationProvisioning")
$ws3.Range("B3").Value = "Ajay Supplier"
$ws3.Range("D2").Value = "SKIP"
$ws3.Range("D3").Value = "PASS"

# --- Sheet "DashboardSetup" ---
$ws4 = $wb.Worksheets.Item("DashboardSetup")
$ws4.Range("B3").Value = "Ajay Supplier"
$ws4.Range("D2").Value = "SKIP"
$ws4.Range("D3").Value = "PASS"

# --- Restore selections to match the authored view state ---
$ws1.Activate()
$ws1.Range("C14").Select() | Out-Null
$ws2.Activate()
$ws2.Range("D3").Select() | Out-Null
$ws3.Activate()
$ws3.Range("C3").Select() | Out-Null
$ws4.Activate()
$ws4.Range("B3").Select() | Out-Null

$ws1.Activate()
